# Weekly refresh of "Fruta / hortaliza" market data: replace the values in
# columns D (Fecha), K (Variedad), L (Calidad), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), Q (Unidad de comercializacion),
# R (Origen), S (Precio $/Kg) and T (Kg / unidad) for every data row (2-25).
# Columns A, B, C, E, F, G, H, I, J (market/product identifiers) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2; D = 44908; K = "Rainier"; L = "Segunda"; M = 250; N = 15000; O = 16000; P = 15600; Q = "$/caja 10 kilos"; R = "Región de O'Higgins"; S = 1560; T = 10 },
    @{ Row = 3; D = 44537; K = "Brooks"; L = "Primera"; M = 200; N = 29000; O = 30000; P = 29500; Q = "$/caja 20 kilos"; R = "Región de O'Higgins"; S = 1475; T = 20 },
    @{ Row = 4; D = 44922; K = "Bing"; L = "Primera"; M = 200; N = 5000; O = 6000; P = 5500; Q = "$/bandeja 10 kilos"; R = "Región del Maule"; S = 550; T = 10 },
    @{ Row = 5; D = 44175; K = "Rainier"; L = "Segunda"; M = 270; N = 25000; O = 26000; P = 25500; Q = "$/caja 18 kilos"; R = "Región de O'Higgins"; S = 1417; T = 18 },
    @{ Row = 6; D = 44571; K = "Brooks"; L = "Segunda"; M = 400; N = 8500; O = 9000; P = 8750; Q = "$/bandeja 10 kilos"; R = "Región de O'Higgins"; S = 875; T = 10 },
    @{ Row = 7; D = 44557; K = "Lapins"; L = "Primera"; M = 250; N = 9000; O = 10000; P = 9500; Q = "$/bandeja 10 kilos"; R = "Provincia de Curicó"; S = 950; T = 10 },
    @{ Row = 8; D = 44901; K = "Bing"; L = "Primera"; M = 500; N = 12000; O = 13000; P = 12500; Q = "$/caja 15 kilos"; R = "Región de O'Higgins"; S = 833; T = 15 },
    @{ Row = 9; D = 44901; K = "Lapins"; L = "Primera"; M = 500; N = 12000; O = 13000; P = 12500; Q = "$/caja 15 kilos"; R = "Región de O'Higgins"; S = 833; T = 15 },
    @{ Row = 10; D = 44568; K = "Santina"; L = "Segunda"; M = 200; N = 15000; O = 16000; P = 15500; Q = "$/bandeja 12 kilos"; R = "Región de O'Higgins"; S = 1292; T = 12 },
    @{ Row = 11; D = 44210; K = "Rainier"; L = "Segunda"; M = 250; N = 21000; O = 22000; P = 21500; Q = "$/caja 18 kilos"; R = "Región de O'Higgins"; S = 1194; T = 18 },
    @{ Row = 12; D = 44229; K = "Santina"; L = "Primera"; M = 250; N = 6500; O = 7000; P = 6750; Q = "$/bandeja 5 kilos"; R = "Provincia de Curicó"; S = 1350; T = 5 },
    @{ Row = 13; D = 44921; K = "Bing"; L = "Primera"; M = 320; N = 7500; O = 8000; P = 7781; Q = "$/bandeja 10 kilos"; R = "Región de O'Higgins"; S = 778; T = 10 },
    @{ Row = 14; D = 44594; K = "Santina"; L = "Primera"; M = 160; N = 5000; O = 6000; P = 5500; Q = "$/bandeja 5 kilos"; R = "Región de O'Higgins"; S = 1100; T = 5 },
    @{ Row = 15; D = 44208; K = "Lapins"; L = "Segunda"; M = 200; N = 10500; O = 11000; P = 10750; Q = "$/bandeja 12 kilos"; R = "Provincia de Curicó"; S = 896; T = 12 },
    @{ Row = 16; D = 44931; K = "Lapins"; L = "Segunda"; M = 250; N = 6000; O = 6500; P = 6250; Q = "$/bandeja 10 kilos"; R = "Región de O'Higgins"; S = 625; T = 10 },
    @{ Row = 17; D = 44931; K = "Lapins"; L = "Segunda"; M = 400; N = 3000; O = 3300; P = 3150; Q = "$/bandeja 5 kilos"; R = "Región de O'Higgins"; S = 630; T = 5 },
    @{ Row = 18; D = 44943; K = "Santina"; L = "Primera"; M = 600; N = 14000; O = 15000; P = 14333; Q = "$/caja 15 kilos"; R = "Región del Maule"; S = 956; T = 15 },
    @{ Row = 19; D = 44532; K = "Brooks"; L = "Primera"; M = 400; N = 27000; O = 28000; P = 27500; Q = "$/bandeja 12 kilos"; R = "Región de O'Higgins"; S = 2292; T = 12 },
    @{ Row = 20; D = 44580; K = "Sweet Heart"; L = "Segunda"; M = 300; N = 7000; O = 8000; P = 7500; Q = "$/bandeja 10 kilos"; R = "Región de O'Higgins"; S = 750; T = 10 },
    @{ Row = 21; D = 44161; K = "Bing"; L = "Primera"; M = 160; N = 39000; O = 40000; P = 39500; Q = "$/caja 20 kilos"; R = "Provincia de Curicó"; S = 1975; T = 20 },
    @{ Row = 22; D = 44917; K = "Bing"; L = "Primera"; M = 400; N = 5000; O = 6000; P = 5625; Q = "$/bandeja 10 kilos"; R = "Región de O'Higgins"; S = 562; T = 10 },
    @{ Row = 23; D = 44917; K = "Santina"; L = "Primera"; M = 400; N = 5000; O = 6000; P = 5500; Q = "$/bandeja 10 kilos"; R = "Región de O'Higgins"; S = 550; T = 10 },
    @{ Row = 24; D = 44914; K = "Brooks"; L = "Primera"; M = 700; N = 7000; O = 8000; P = 7429; Q = "$/bandeja 10 kilos"; R = "Región de O'Higgins"; S = 743; T = 10 },
    @{ Row = 25; D = 44914; K = "Lapins"; L = "Primera"; M = 550; N = 7000; O = 8000; P = 7455; Q = "$/bandeja 10 kilos"; R = "Región de O'Higgins"; S = 746; T = 10 }
)

foreach ($row in $rows) {
    $r = $row.Row
    $ws.Range("D$r").Value = $row.D
    $ws.Range("K$r").Value = $row.K
    $ws.Range("L$r").Value = $row.L
    $ws.Range("M$r").Value = $row.M
    $ws.Range("N$r").Value = $row.N
    $ws.Range("O$r").Value = $row.O
    $ws.Range("P$r").Value = $row.P
    $ws.Range("Q$r").Value = $row.Q
    $ws.Range("R$r").Value = $row.R
    $ws.Range("S$r").Value = $row.S
    $ws.Range("T$r").Value = $row.T
}
